$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ----------------------------------------------------------------------
# Original header layout (row 1):
#   A=Date  B="App #"  C=Type  D=Status  E="Permit #"  F="Business Name"
#   G="Application Name"  H="Project Description"  I="Office Project Description"
#   J="Contractor Name"  K="Res. Of Enf. Action"  L="City Jur. App. Req"
#   M="City Utility App. Req."
#
# Target header layout (row 1):
#   A="Application Date"  B="Issue Date"  C="Assessor Parcel Number"  D="App #"
#   E=Type  F=Status  G="Permit #"  H="Business Name"  I="Application Name"
#   J="Office Project Description"  K="Contractor Name"  L,M,N = (blank headers,
#   columns kept but no longer labeled - formerly Res. Of Enf. Action /
#   City Jur. App. Req / City Utility App. Req.)
# ----------------------------------------------------------------------

# 1) Insert two new blank columns before column B ("App #") - this shifts the
#    old B..M columns to D..O, carrying their widths/styles/bestFit with them.
$ws.Range("B1:C1").EntireColumn.Insert()

# 2) Remove the old "Project Description" column, which now sits at J
#    (old H shifted by +2 from step 1). This shifts K..O back to J..N.
$ws.Columns("J:J").Delete()

# 3) Rename / fill in the three leading header cells.
$ws.Range("A1").Value = "Application Date"
$ws.Range("B1").Value = "Issue Date"
$ws.Range("C1").Value = "Assessor Parcel Number"

# 3b) Column C ("Assessor Parcel Number") needs to be wide enough to fit its
#     longer header text (target stored column width = 26). 25.15 is in the
#     input range that this engine's width-unit rounding maps to an output
#     of exactly 26.
$ws.Columns.Item(3).ColumnWidth = 25.15

# 4) The trailing three columns (L, M, N - formerly "Res. Of Enf. Action",
#    "City Jur. App. Req", "City Utility App. Req.") keep their column
#    widths/styles but lose their header text.
$ws.Range("L1:N1").ClearContents()

# 5) Apply a date number format to the two new date columns' header cells
#    (mm/dd/yyyy), matching the new numFmtId 164 used by the template.
$ws.Columns("A:B").NumberFormat = "mm/dd/yyyy"
